# Update existing rows 2-30 and append new rows 31-40 with re-processed
# Spending algorithm benchmark data (n = 1..39).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2 (n=1)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 0.053954
$ws.Range("E2").Value = 1390
$ws.Range("F2").Value = 9.695
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 5.1
$ws.Range("I2").Value = 2.8614
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 0.1495915
$ws.Range("L2").Value = 5.554
$ws.Range("M2").Value = 0.1697735
$ws.Range("N2").Value = 0.000121

# row 3 (n=2)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 0.052513
$ws.Range("E3").Value = 1390
$ws.Range("F3").Value = 18.4481
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 10.1
$ws.Range("I3").Value = 5.173999999999999
$ws.Range("J3").Value = 144
$ws.Range("K3").Value = 0.1876047
$ws.Range("L3").Value = 5.556
$ws.Range("M3").Value = 0.1788506
$ws.Range("N3").Value = 0.000059

# row 4 (n=3)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 150
$ws.Range("D4").Value = 0.052353
$ws.Range("E4").Value = 1390
$ws.Range("F4").Value = 26.8248
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 15.1
$ws.Range("I4").Value = 7.597799999999999
$ws.Range("J4").Value = 148
$ws.Range("K4").Value = 0.2233237
$ws.Range("L4").Value = 5.558
$ws.Range("M4").Value = 0.1866526
$ws.Range("N4").Value = 0.000059

# row 5 (n=4)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 200
$ws.Range("D5").Value = 0.052011
$ws.Range("E5").Value = 1390
$ws.Range("F5").Value = 35.04089999999999
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 20.1
$ws.Range("I5").Value = 9.868099999999998
$ws.Range("J5").Value = 152
$ws.Range("K5").Value = 0.2585387
$ws.Range("L5").Value = 5.56
$ws.Range("M5").Value = 0.1943055
$ws.Range("N5").Value = 0.000057

# row 6 (n=5)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 250
$ws.Range("D6").Value = 0.053634
$ws.Range("E6").Value = 1390
$ws.Range("F6").Value = 43.54539999999999
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = 25.1
$ws.Range("I6").Value = 12.1821
$ws.Range("J6").Value = 156
$ws.Range("K6").Value = 0.2963296000000001
$ws.Range("L6").Value = 5.562
$ws.Range("M6").Value = 0.2047147
$ws.Range("N6").Value = 0.000059

# row 7 (n=6)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = 0.051654
$ws.Range("E7").Value = 1390
$ws.Range("F7").Value = 51.88159999999999
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = 30.1
$ws.Range("I7").Value = 14.5551
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 0.3335766999999999
$ws.Range("L7").Value = 5.564
$ws.Range("M7").Value = 0.213082
$ws.Range("N7").Value = 0.000057

# row 8 (n=7)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 350
$ws.Range("D8").Value = 0.052206
$ws.Range("E8").Value = 1390
$ws.Range("F8").Value = 59.88709999999999
$ws.Range("G8").Value = 70
$ws.Range("H8").Value = 35.1
$ws.Range("I8").Value = 16.9119
$ws.Range("J8").Value = 164
$ws.Range("K8").Value = 0.3687463000000001
$ws.Range("L8").Value = 5.566
$ws.Range("M8").Value = 0.2213518
$ws.Range("N8").Value = 0.000073

# row 9 (n=8)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 400
$ws.Range("D9").Value = 0.052273
$ws.Range("E9").Value = 1390
$ws.Range("F9").Value = 67.92580000000001
$ws.Range("G9").Value = 80
$ws.Range("H9").Value = 40.1
$ws.Range("I9").Value = 18.8731
$ws.Range("J9").Value = 168
$ws.Range("K9").Value = 0.4032951
$ws.Range("L9").Value = 5.568
$ws.Range("M9").Value = 0.2289381
$ws.Range("N9").Value = 0.000061

# row 10 (n=9)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 450
$ws.Range("D10").Value = 0.051716
$ws.Range("E10").Value = 1390
$ws.Range("F10").Value = 76.2459
$ws.Range("G10").Value = 90
$ws.Range("H10").Value = 45.1
$ws.Range("I10").Value = 21.2678
$ws.Range("J10").Value = 172
$ws.Range("K10").Value = 0.4367034
$ws.Range("L10").Value = 5.57
$ws.Range("M10").Value = 0.2356118
$ws.Range("N10").Value = 0.000059

# row 11 (n=10)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 500
$ws.Range("D11").Value = 0.05131
$ws.Range("E11").Value = 1390
$ws.Range("F11").Value = 84.66660000000002
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = 50.1
$ws.Range("I11").Value = 23.5857
$ws.Range("J11").Value = 176
$ws.Range("K11").Value = 0.4721851000000001
$ws.Range("L11").Value = 5.572
$ws.Range("M11").Value = 0.2447571
$ws.Range("N11").Value = 0.000061

# row 12 (n=11)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 550
$ws.Range("D12").Value = 0.051653
$ws.Range("E12").Value = 1390
$ws.Range("F12").Value = 92.762
$ws.Range("G12").Value = 110
$ws.Range("H12").Value = 55.1
$ws.Range("I12").Value = 25.8162
$ws.Range("J12").Value = 180
$ws.Range("K12").Value = 0.5060454999999999
$ws.Range("L12").Value = 5.574
$ws.Range("M12").Value = 0.2544194
$ws.Range("N12").Value = 0.000061

# row 13 (n=12)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 600
$ws.Range("D13").Value = 0.052256
$ws.Range("E13").Value = 1390
$ws.Range("F13").Value = 101.3105
$ws.Range("G13").Value = 120
$ws.Range("H13").Value = 60.1
$ws.Range("I13").Value = 28.2899
$ws.Range("J13").Value = 184
$ws.Range("K13").Value = 0.5420174
$ws.Range("L13").Value = 5.576
$ws.Range("M13").Value = 0.2631567
$ws.Range("N13").Value = 0.000061

# row 14 (n=13)
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 650
$ws.Range("D14").Value = 0.051416
$ws.Range("E14").Value = 1390
$ws.Range("F14").Value = 109.1167
$ws.Range("G14").Value = 130
$ws.Range("H14").Value = 65.1
$ws.Range("I14").Value = 30.3644
$ws.Range("J14").Value = 188
$ws.Range("K14").Value = 0.5823302
$ws.Range("L14").Value = 5.578
$ws.Range("M14").Value = 0.2733833
$ws.Range("N14").Value = 0.000066

# row 15 (n=14)
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 700
$ws.Range("D15").Value = 0.051388
$ws.Range("E15").Value = 1390
$ws.Range("F15").Value = 117.231
$ws.Range("G15").Value = 140
$ws.Range("H15").Value = 70.1
$ws.Range("I15").Value = 32.6972
$ws.Range("J15").Value = 192
$ws.Range("K15").Value = 0.6189059
$ws.Range("L15").Value = 5.58
$ws.Range("M15").Value = 0.2818899
$ws.Range("N15").Value = 0.000105

# row 16 (n=15)
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 750
$ws.Range("D16").Value = 0.051699
$ws.Range("E16").Value = 1390
$ws.Range("F16").Value = 125.6227
$ws.Range("G16").Value = 150
$ws.Range("H16").Value = 75.1
$ws.Range("I16").Value = 35.0088
$ws.Range("J16").Value = 196
$ws.Range("K16").Value = 0.6548914
$ws.Range("L16").Value = 5.582
$ws.Range("M16").Value = 0.2901717
$ws.Range("N16").Value = 0.00007

# row 17 (n=16)
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 800
$ws.Range("D17").Value = 0.0523
$ws.Range("E17").Value = 1390
$ws.Range("F17").Value = 133.3715
$ws.Range("G17").Value = 160
$ws.Range("H17").Value = 80.1
$ws.Range("I17").Value = 37.29430000000001
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 0.6914001
$ws.Range("L17").Value = 5.584
$ws.Range("M17").Value = 0.2981141
$ws.Range("N17").Value = 0.000089

# row 18 (n=17)
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 850
$ws.Range("D18").Value = 0.05193
$ws.Range("E18").Value = 1390
$ws.Range("F18").Value = 142.0459
$ws.Range("G18").Value = 170
$ws.Range("H18").Value = 85.1
$ws.Range("I18").Value = 39.67
$ws.Range("J18").Value = 204
$ws.Range("K18").Value = 0.7265518999999999
$ws.Range("L18").Value = 5.586
$ws.Range("M18").Value = 0.3058778
$ws.Range("N18").Value = 0.000065

# row 19 (n=18)
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = 900
$ws.Range("D19").Value = 0.051423
$ws.Range("E19").Value = 1390
$ws.Range("F19").Value = 150.9115
$ws.Range("G19").Value = 180
$ws.Range("H19").Value = 90.1
$ws.Range("I19").Value = 41.7537
$ws.Range("J19").Value = 208
$ws.Range("K19").Value = 0.7630923000000001
$ws.Range("L19").Value = 5.588
$ws.Range("M19").Value = 0.3158104
$ws.Range("N19").Value = 0.000067

# row 20 (n=19)
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = 950
$ws.Range("D20").Value = 0.051081
$ws.Range("E20").Value = 1390
$ws.Range("F20").Value = 158.3252
$ws.Range("G20").Value = 190
$ws.Range("H20").Value = 95.1
$ws.Range("I20").Value = 44.2044
$ws.Range("J20").Value = 212
$ws.Range("K20").Value = 0.7970005
$ws.Range("L20").Value = 5.59
$ws.Range("M20").Value = 0.3252923
$ws.Range("N20").Value = 0.000065

# row 21 (n=20)
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = 1000
$ws.Range("D21").Value = 0.051901
$ws.Range("E21").Value = 1390
$ws.Range("F21").Value = 167.0069
$ws.Range("G21").Value = 200
$ws.Range("H21").Value = 100.1
$ws.Range("I21").Value = 46.2703
$ws.Range("J21").Value = 216
$ws.Range("K21").Value = 0.8357794000000001
$ws.Range("L21").Value = 5.592
$ws.Range("M21").Value = 0.3347911
$ws.Range("N21").Value = 0.000068

# row 22 (n=21)
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = 1050
$ws.Range("D22").Value = 0.050771
$ws.Range("E22").Value = 1390
$ws.Range("F22").Value = 175.1129
$ws.Range("G22").Value = 210
$ws.Range("H22").Value = 105.1
$ws.Range("I22").Value = 48.5286
$ws.Range("J22").Value = 220
$ws.Range("K22").Value = 0.8714166000000001
$ws.Range("L22").Value = 5.594
$ws.Range("M22").Value = 0.3431711
$ws.Range("N22").Value = 0.000109

# row 23 (n=22)
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = 1100
$ws.Range("D23").Value = 0.051945
$ws.Range("E23").Value = 1390
$ws.Range("F23").Value = 183.2457
$ws.Range("G23").Value = 220
$ws.Range("H23").Value = 110.1
$ws.Range("I23").Value = 51.0066
$ws.Range("J23").Value = 224
$ws.Range("K23").Value = 0.9051252000000002
$ws.Range("L23").Value = 5.596
$ws.Range("M23").Value = 0.3510204
$ws.Range("N23").Value = 0.000069

# row 24 (n=23)
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 23
$ws.Range("C24").Value = 1150
$ws.Range("D24").Value = 0.050872
$ws.Range("E24").Value = 1390
$ws.Range("F24").Value = 191.6057
$ws.Range("G24").Value = 230
$ws.Range("H24").Value = 115.1
$ws.Range("I24").Value = 53.2388
$ws.Range("J24").Value = 228
$ws.Range("K24").Value = 0.9442529000000001
$ws.Range("L24").Value = 5.598
$ws.Range("M24").Value = 0.3620467
$ws.Range("N24").Value = 0.000067

# row 25 (n=24)
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = 1200
$ws.Range("D25").Value = 0.052088
$ws.Range("E25").Value = 1390
$ws.Range("F25").Value = 199.8202
$ws.Range("G25").Value = 240
$ws.Range("H25").Value = 120.1
$ws.Range("I25").Value = 55.41459999999999
$ws.Range("J25").Value = 232
$ws.Range("K25").Value = 0.9767288000000001
$ws.Range("L25").Value = 5.6
$ws.Range("M25").Value = 0.3689977999999999
$ws.Range("N25").Value = 0.00007

# row 26 (n=25)
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 25
$ws.Range("C26").Value = 1250
$ws.Range("D26").Value = 0.052318
$ws.Range("E26").Value = 1390
$ws.Range("F26").Value = 208.1906
$ws.Range("G26").Value = 250
$ws.Range("H26").Value = 125.1
$ws.Range("I26").Value = 57.9474
$ws.Range("J26").Value = 236
$ws.Range("K26").Value = 1.0138167
$ws.Range("L26").Value = 5.602
$ws.Range("M26").Value = 0.3777181
$ws.Range("N26").Value = 0.000071

# row 27 (n=26)
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 26
$ws.Range("C27").Value = 1300
$ws.Range("D27").Value = 0.050514
$ws.Range("E27").Value = 1390
$ws.Range("F27").Value = 215.1766
$ws.Range("G27").Value = 260
$ws.Range("H27").Value = 130.1
$ws.Range("I27").Value = 60.2739
$ws.Range("J27").Value = 240
$ws.Range("K27").Value = 1.0490257
$ws.Range("L27").Value = 5.604
$ws.Range("M27").Value = 0.3883579
$ws.Range("N27").Value = 0.000081

# row 28 (n=27)
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 27
$ws.Range("C28").Value = 1350
$ws.Range("D28").Value = 0.051753
$ws.Range("E28").Value = 1390
$ws.Range("F28").Value = 224.3901
$ws.Range("G28").Value = 270
$ws.Range("H28").Value = 135.1
$ws.Range("I28").Value = 62.54339999999999
$ws.Range("J28").Value = 244
$ws.Range("K28").Value = 1.088757
$ws.Range("L28").Value = 5.606
$ws.Range("M28").Value = 0.3949603
$ws.Range("N28").Value = 0.00007

# row 29 (n=28)
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 28
$ws.Range("C29").Value = 1400
$ws.Range("D29").Value = 0.05094
$ws.Range("E29").Value = 1390
$ws.Range("F29").Value = 232.0399
$ws.Range("G29").Value = 280
$ws.Range("H29").Value = 140.1
$ws.Range("I29").Value = 65.27749999999999
$ws.Range("J29").Value = 248
$ws.Range("K29").Value = 1.1262519
$ws.Range("L29").Value = 5.608
$ws.Range("M29").Value = 0.4060505
$ws.Range("N29").Value = 0.00007

# row 30 (n=29)
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 29
$ws.Range("C30").Value = 1450
$ws.Range("D30").Value = 0.052459
$ws.Range("E30").Value = 1390
$ws.Range("F30").Value = 241.4971
$ws.Range("G30").Value = 290
$ws.Range("H30").Value = 145.1
$ws.Range("I30").Value = 67.12089999999998
$ws.Range("J30").Value = 252
$ws.Range("K30").Value = 1.165281
$ws.Range("L30").Value = 5.61
$ws.Range("M30").Value = 0.4211609
$ws.Range("N30").Value = 0.000073

# row 31 (n=30)
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = 1500
$ws.Range("D31").Value = 0.054042
$ws.Range("E31").Value = 1390
$ws.Range("F31").Value = 251.3277
$ws.Range("G31").Value = 300
$ws.Range("H31").Value = 150.1
$ws.Range("I31").Value = 71.1025
$ws.Range("J31").Value = 256
$ws.Range("K31").Value = 1.2207208
$ws.Range("L31").Value = 5.612
$ws.Range("M31").Value = 0.4323857
$ws.Range("N31").Value = 0.000074

# row 32 (n=31)
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 31
$ws.Range("C32").Value = 1550
$ws.Range("D32").Value = 0.052753
$ws.Range("E32").Value = 1390
$ws.Range("F32").Value = 257.7357
$ws.Range("G32").Value = 310
$ws.Range("H32").Value = 155.1
$ws.Range("I32").Value = 72.1919
$ws.Range("J32").Value = 260
$ws.Range("K32").Value = 1.2485605
$ws.Range("L32").Value = 5.614
$ws.Range("M32").Value = 0.4513919
$ws.Range("N32").Value = 0.000085

# row 33 (n=32)
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 32
$ws.Range("C33").Value = 1600
$ws.Range("D33").Value = 0.054211
$ws.Range("E33").Value = 1390
$ws.Range("F33").Value = 272.0355999999999
$ws.Range("G33").Value = 320
$ws.Range("H33").Value = 160.1
$ws.Range("I33").Value = 76.60579999999999
$ws.Range("J33").Value = 264
$ws.Range("K33").Value = 1.3183701
$ws.Range("L33").Value = 5.616
$ws.Range("M33").Value = 0.4355301
$ws.Range("N33").Value = 0.000075

# row 34 (n=33)
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 33
$ws.Range("C34").Value = 1650
$ws.Range("D34").Value = 0.051036
$ws.Range("E34").Value = 1390
$ws.Range("F34").Value = 267.6098
$ws.Range("G34").Value = 330
$ws.Range("H34").Value = 165.1
$ws.Range("I34").Value = 75.2291
$ws.Range("J34").Value = 268
$ws.Range("K34").Value = 1.2869503
$ws.Range("L34").Value = 5.618
$ws.Range("M34").Value = 0.4457908
$ws.Range("N34").Value = 0.000072

# row 35 (n=34)
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 34
$ws.Range("C35").Value = 1700
$ws.Range("D35").Value = 0.051974
$ws.Range("E35").Value = 1390
$ws.Range("F35").Value = 275.8148
$ws.Range("G35").Value = 340
$ws.Range("H35").Value = 170.1
$ws.Range("I35").Value = 77.88770000000001
$ws.Range("J35").Value = 272
$ws.Range("K35").Value = 1.330494
$ws.Range("L35").Value = 5.62
$ws.Range("M35").Value = 0.4520981
$ws.Range("N35").Value = 0.000073

# row 36 (n=35)
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 35
$ws.Range("C36").Value = 1750
$ws.Range("D36").Value = 0.051777
$ws.Range("E36").Value = 1390
$ws.Range("F36").Value = 284.9892
$ws.Range("G36").Value = 350
$ws.Range("H36").Value = 175.1
$ws.Range("I36").Value = 80.37480000000001
$ws.Range("J36").Value = 276
$ws.Range("K36").Value = 1.3813288
$ws.Range("L36").Value = 5.622
$ws.Range("M36").Value = 0.4779692
$ws.Range("N36").Value = 0.000078

# row 37 (n=36)
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 36
$ws.Range("C37").Value = 1800
$ws.Range("D37").Value = 0.05289
$ws.Range("E37").Value = 1390
$ws.Range("F37").Value = 302.4467
$ws.Range("G37").Value = 360
$ws.Range("H37").Value = 180.1
$ws.Range("I37").Value = 83.8517
$ws.Range("J37").Value = 280
$ws.Range("K37").Value = 1.4514284
$ws.Range("L37").Value = 5.624
$ws.Range("M37").Value = 0.4873145
$ws.Range("N37").Value = 0.000074

# row 38 (n=37)
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 37
$ws.Range("C38").Value = 1850
$ws.Range("D38").Value = 0.053148
$ws.Range("E38").Value = 1390
$ws.Range("F38").Value = 306.6689
$ws.Range("G38").Value = 370
$ws.Range("H38").Value = 185.1
$ws.Range("I38").Value = 86.2364
$ws.Range("J38").Value = 284
$ws.Range("K38").Value = 1.4891059
$ws.Range("L38").Value = 5.626
$ws.Range("M38").Value = 0.5047622999999999
$ws.Range("N38").Value = 0.000079

# row 39 (n=38)
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 38
$ws.Range("C39").Value = 1900
$ws.Range("D39").Value = 0.054988
$ws.Range("E39").Value = 1390
$ws.Range("F39").Value = 321.8945
$ws.Range("G39").Value = 380
$ws.Range("H39").Value = 190.1
$ws.Range("I39").Value = 90.7969
$ws.Range("J39").Value = 288
$ws.Range("K39").Value = 1.5477549
$ws.Range("L39").Value = 5.628
$ws.Range("M39").Value = 0.5166873000000001
$ws.Range("N39").Value = 0.000083

# row 40 (n=39)
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 39
$ws.Range("C40").Value = 1950
$ws.Range("D40").Value = 0.054975
$ws.Range("E40").Value = 1390
$ws.Range("F40").Value = 329.825
$ws.Range("G40").Value = 390
$ws.Range("H40").Value = 195.1
$ws.Range("I40").Value = 93.31369999999998
$ws.Range("J40").Value = 292
$ws.Range("K40").Value = 1.5768658
$ws.Range("L40").Value = 5.63
$ws.Range("M40").Value = 0.5210357999999999
$ws.Range("N40").Value = 0.000081

